$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5504.759
$ws.Range("I86").Value = 5554.8887
$ws.Range("J86").Value = 5422.727
$ws.Range("K86").Value = 5554.8887
$ws.Range("L86").Value = 5422.727
$ws.Range("M86").Value = -4431.8887
$ws.Range("N86").Value = -7668.727

$ws.Range("H89").Value = 5504.759
$ws.Range("I89").Value = 5554.8887
$ws.Range("J89").Value = 5422.727
$ws.Range("K89").Value = 27774.4435
$ws.Range("L89").Value = 27113.635
$ws.Range("M89").Value = -22158.4435
$ws.Range("N89").Value = -38345.63499999999

$ws.Range("H113").Value = 3402.5
$ws.Range("I113").Value = 3402.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3402.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -148.5
$ws.Range("N113").ClearContents()

$ws.Range("H121").Value = 1685223.4
$ws.Range("J121").Value = 1895432.9
$ws.Range("L121").Value = 5686298.699999999
$ws.Range("N121").Value = -5689792.699999999

$ws.Range("H131").Value = 6605.9375
$ws.Range("I131").Value = 6369.6
$ws.Range("J131").Value = 6999.8335
$ws.Range("K131").Value = 19108.8
$ws.Range("L131").Value = 20999.5005
$ws.Range("M131").Value = -14068.8
$ws.Range("N131").Value = -31079.5005

$ws.Range("H137").Value = 11442.137
$ws.Range("J137").Value = 14979.862
$ws.Range("L137").Value = 44939.586
$ws.Range("N137").Value = -50039.586

$ws.Range("H138").Value = 6464.8
$ws.Range("J138").Value = 6079.4707
$ws.Range("L138").Value = 18238.4121
$ws.Range("N138").Value = -28518.4121

$ws.Range("H141").Value = 2305.8572
$ws.Range("I141").Value = 2473.5
$ws.Range("J141").Value = 2082.3333
$ws.Range("K141").Value = 7420.5
$ws.Range("L141").Value = 6246.999899999999
$ws.Range("M141").Value = -2240.5
$ws.Range("N141").Value = -16606.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5007129.5
$ws.Range("I32").Value = 5380863
$ws.Range("K32").Value = 5380863
$ws.Range("M32").Value = -5380576

$ws.Range("H74").Value = 5213624.5
$ws.Range("I74").Value = 7354530.5
$ws.Range("J74").Value = 14280.857
$ws.Range("K74").Value = 7354530.5
$ws.Range("L74").Value = 14280.857
$ws.Range("M74").Value = -7353656.5
$ws.Range("N74").Value = -16028.857

$ws.Range("H77").Value = 5213624.5
$ws.Range("I77").Value = 7354530.5
$ws.Range("J77").Value = 14280.857
$ws.Range("K77").Value = 36772652.5
$ws.Range("L77").Value = 71404.285
$ws.Range("M77").Value = -36768284.5
$ws.Range("N77").Value = -80140.285

$ws.Range("H124").Value = 80231.8
$ws.Range("J124").Value = 80231.8
$ws.Range("L124").Value = 80231.8
$ws.Range("N124").Value = -90051.8

$ws.Range("H125").Value = 60675.75
$ws.Range("J125").Value = 60675.75
$ws.Range("L125").Value = 60675.75
$ws.Range("N125").Value = -70515.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 29784.555
$ws.Range("I96").Value = 11899.333
$ws.Range("J96").Value = 65555
$ws.Range("K96").Value = 11899.333
$ws.Range("L96").Value = 65555
$ws.Range("M96").Value = -9153.333000000001
$ws.Range("N96").Value = -71047

$ws.Range("H107").Value = 2473.15
$ws.Range("I107").Value = 2172.3572
$ws.Range("J107").Value = 3175
$ws.Range("K107").Value = 2172.3572
$ws.Range("L107").Value = 3175
$ws.Range("M107").Value = -252.3571999999999
$ws.Range("N107").Value = -7015

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 513364.25
$ws.Range("I31").Value = 13296.643
$ws.Range("J31").Value = 732143.9
$ws.Range("K31").Value = 13296.643
$ws.Range("L31").Value = 732143.9
$ws.Range("M31").Value = -13001.643
$ws.Range("N31").Value = -732733.9

$ws.Range("H34").Value = 513364.25
$ws.Range("I34").Value = 13296.643
$ws.Range("J34").Value = 732143.9
$ws.Range("K34").Value = 13296.643
$ws.Range("L34").Value = 732143.9
$ws.Range("M34").Value = -13094.643
$ws.Range("N34").Value = -732547.9

$ws.Range("H43").Value = 117000
$ws.Range("J43").Value = 117000
$ws.Range("L43").Value = 117000
$ws.Range("N43").Value = -117368

$ws.Range("H96").Value = 17555
$ws.Range("J96").Value = 21444
$ws.Range("L96").Value = 21444
$ws.Range("N96").Value = -26936

$ws.Range("H101").Value = 117000
$ws.Range("J101").Value = 117000
$ws.Range("L101").Value = 117000
$ws.Range("N101").Value = -123490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 1999.6666
$ws.Range("I35").Value = 1999.6666
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 5998.9998
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -5710.9998
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18092646
$ws.Range("J11").Value = 33834.5
$ws.Range("L11").Value = 33834.5
$ws.Range("N11").Value = -34112.5

$ws.Range("H18").Value = 37666
$ws.Range("I18").Value = 37666
$ws.Range("K18").Value = 37666
$ws.Range("M18").Value = -37373

$ws.Range("H43").Value = 20583.2
$ws.Range("I43").Value = 15979.25
$ws.Range("J43").Value = 38999
$ws.Range("K43").Value = 15979.25
$ws.Range("L43").Value = 38999
$ws.Range("M43").Value = -15828.25
$ws.Range("N43").Value = -39301

$ws.Range("H95").Value = 57517.8
$ws.Range("J95").Value = 57517.8
$ws.Range("L95").Value = 57517.8
$ws.Range("N95").Value = -63009.8

$ws.Range("H126").Value = 3807.25
$ws.Range("I126").Value = 3624.625
$ws.Range("J126").Value = 4172.5
$ws.Range("K126").Value = 10873.875
$ws.Range("L126").Value = 12517.5
$ws.Range("M126").Value = -8403.875
$ws.Range("N126").Value = -17457.5

$ws.Range("H130").Value = 59644.285
$ws.Range("J130").Value = 59644.285
$ws.Range("L130").Value = 59644.285
$ws.Range("N130").Value = -69684.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45604.543
$ws.Range("I7").Value = 3052.0625
$ws.Range("J7").Value = 130709.5
$ws.Range("K7").Value = 3052.0625
$ws.Range("L7").Value = 130709.5
$ws.Range("M7").Value = -2940.0625
$ws.Range("N7").Value = -130933.5

$ws.Range("H40").Value = 3831.1292
$ws.Range("I40").Value = 3115.111
$ws.Range("J40").Value = 4822.5386
$ws.Range("K40").Value = 3115.111
$ws.Range("L40").Value = 4822.5386
$ws.Range("M40").Value = -2979.111
$ws.Range("N40").Value = -5094.5386

$ws.Range("H82").Value = 1282.6316
$ws.Range("I82").Value = 890.1818
$ws.Range("J82").Value = 1822.25
$ws.Range("K82").Value = 890.1818
$ws.Range("L82").Value = 1822.25
$ws.Range("M82").Value = -529.1818
$ws.Range("N82").Value = -2544.25

$ws.Range("H85").Value = 1282.6316
$ws.Range("I85").Value = 890.1818
$ws.Range("J85").Value = 1822.25
$ws.Range("K85").Value = 890.1818
$ws.Range("L85").Value = 1822.25
$ws.Range("M85").Value = 357.8182
$ws.Range("N85").Value = -4318.25

$ws.Range("H92").Value = 103956
$ws.Range("J92").Value = 103956
$ws.Range("L92").Value = 103956
$ws.Range("N92").Value = -108948

$ws.Range("H96").Value = 64065.668
$ws.Range("I96").Value = 5000
$ws.Range("K96").Value = 5000
$ws.Range("M96").Value = -2254

$ws.Range("H126").Value = 45604.543
$ws.Range("I126").Value = 3052.0625
$ws.Range("J126").Value = 130709.5
$ws.Range("K126").Value = 9156.1875
$ws.Range("L126").Value = 392128.5
$ws.Range("M126").Value = -6686.1875
$ws.Range("N126").Value = -397068.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 109999.5
$ws.Range("I81").Value = 199999
$ws.Range("J81").Value = 20000
$ws.Range("K81").Value = 399998
$ws.Range("L81").Value = 40000
$ws.Range("M81").Value = -398937
$ws.Range("N81").Value = -42122

$ws.Range("H84").Value = 109999.5
$ws.Range("I84").Value = 199999
$ws.Range("J84").Value = 20000
$ws.Range("K84").Value = 1999990
$ws.Range("L84").Value = 200000
$ws.Range("M84").Value = -1994686
$ws.Range("N84").Value = -210608
